$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 799.53845
$ws.Range("I28").Value = 832.8333
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 832.8333
$ws.Range("L28").Value = 400
$ws.Range("M28").Value = -347.8333

$ws.Range("H46").Value = 20819.45
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 23088.277
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 69264.83099999999
$ws.Range("M46").Value = -1081
$ws.Range("N46").Value = -69502.83099999999

$ws.Range("H60").Value = 20819.45
$ws.Range("I60").Value = 400
$ws.Range("J60").Value = 23088.277
$ws.Range("K60").Value = 1200
$ws.Range("L60").Value = 69264.83099999999
$ws.Range("M60").Value = -716
$ws.Range("N60").Value = -70232.83099999999

$ws.Range("H62").Value = 3455.95
$ws.Range("I62").Value = 3068.4375
$ws.Range("J62").Value = 5006
$ws.Range("K62").Value = 3068.4375
$ws.Range("L62").Value = 5006
$ws.Range("M62").Value = -2444.4375

$ws.Range("H65").Value = 3455.95
$ws.Range("I65").Value = 3068.4375
$ws.Range("J65").Value = 5006
$ws.Range("K65").Value = 15342.1875
$ws.Range("L65").Value = 25030
$ws.Range("M65").Value = -12222.1875

$ws.Range("H100").Value = 11942776
$ws.Range("I100").Value = 16668101
$ws.Range("J100").Value = 129463.875
$ws.Range("K100").Value = 16668101
$ws.Range("L100").Value = 129463.875
$ws.Range("M100").Value = -16667560
$ws.Range("N100").Value = -130545.875

$ws.Range("H106").Value = 2000
$ws.Range("I106").Value = 1500
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1500
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -869

$ws.Range("H129").Value = 1091.3914
$ws.Range("I129").Value = 632.3333
$ws.Range("J129").Value = 1123.4186
$ws.Range("K129").Value = 1896.9999
$ws.Range("L129").Value = 3370.2558
$ws.Range("M129").Value = 3103.0001
$ws.Range("N129").Value = -13370.2558

$ws.Range("H138").Value = 2905.013
$ws.Range("I138").Value = 1473.18
$ws.Range("J138").Value = 5461.857
$ws.Range("K138").Value = 4419.54
$ws.Range("L138").Value = 16385.571
$ws.Range("M138").Value = 720.46
$ws.Range("N138").Value = -26665.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 50000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 50000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 50000
$ws.Range("N9").Value = -50340

$ws.Range("H20").Value = 50000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 50000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 50000
$ws.Range("N20").Value = -50540

$ws.Range("H61").Value = 1165.8
$ws.Range("I61").Value = 676.3514
$ws.Range("J61").Value = 2171.889
$ws.Range("K61").Value = 676.3514
$ws.Range("L61").Value = 2171.889
$ws.Range("M61").Value = -464.3514
$ws.Range("N61").Value = -2595.889

$ws.Range("H88").Value = 2400
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2400
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2400
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3212

$ws.Range("H91").Value = 2400
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2400
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2400
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5208

$ws.Range("H132").Value = 2363.5151
$ws.Range("I132").Value = 1242.5714
$ws.Range("J132").Value = 3189.4736
$ws.Range("K132").Value = 3727.7142
$ws.Range("L132").Value = 9568.4208
$ws.Range("M132").Value = -1197.7142
$ws.Range("N132").Value = -14628.4208

$ws.Range("H134").Value = 36209.855
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 36209.855
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 36209.855
$ws.Range("N134").Value = -46349.855

$ws.Range("H136").Value = 1165.8
$ws.Range("I136").Value = 676.3514
$ws.Range("J136").Value = 2171.889
$ws.Range("K136").Value = 2029.0542
$ws.Range("L136").Value = 6515.667
$ws.Range("M136").Value = 520.9458
$ws.Range("N136").Value = -11615.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2678.4119
$ws.Range("I86").Value = 3006.9092
$ws.Range("J86").Value = 2076.1667
$ws.Range("K86").Value = 3006.9092
$ws.Range("L86").Value = 2076.1667
$ws.Range("M86").Value = -1883.9092
$ws.Range("N86").Value = -4322.1667

$ws.Range("H89").Value = 2678.4119
$ws.Range("I89").Value = 3006.9092
$ws.Range("J89").Value = 2076.1667
$ws.Range("K89").Value = 15034.546
$ws.Range("L89").Value = 10380.8335
$ws.Range("M89").Value = -9418.546
$ws.Range("N89").Value = -21612.8335

$ws.Range("H107").Value = 522
$ws.Range("I107").Value = 405.65
$ws.Range("J107").Value = 1297.6666
$ws.Range("K107").Value = 405.65
$ws.Range("L107").Value = 1297.6666
$ws.Range("M107").Value = 1514.35

$ws.Range("H134").Value = 1478.4082
$ws.Range("I134").Value = 998.75
$ws.Range("J134").Value = 2381.2942
$ws.Range("K134").Value = 2996.25
$ws.Range("L134").Value = 7143.882599999999
$ws.Range("M134").Value = -461.25
$ws.Range("N134").Value = -12213.8826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2276.5
$ws.Range("I105").Value = 2650.9092
$ws.Range("J105").Value = 903.6667
$ws.Range("K105").Value = 2650.9092
$ws.Range("L105").Value = 903.6667
$ws.Range("M105").Value = -903.9092000000001
$ws.Range("N105").Value = -4397.6667

$ws.Range("H132").Value = 1268.6
$ws.Range("I132").Value = 601.7805
$ws.Range("J132").Value = 3221.4285
$ws.Range("K132").Value = 1805.3415
$ws.Range("L132").Value = 9664.2855
$ws.Range("M132").Value = 724.6585
$ws.Range("N132").Value = -14724.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 260564
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 260564
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 781692
$ws.Range("N76").Value = -782458

$ws.Range("H79").Value = 260564
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 260564
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 781692
$ws.Range("N79").Value = -784344

$ws.Range("H109").Value = 2765.4614
$ws.Range("I109").Value = 521.4
$ws.Range("J109").Value = 3299.762
$ws.Range("K109").Value = 1564.2
$ws.Range("L109").Value = 9899.286
$ws.Range("M109").Value = -524.1999999999998
$ws.Range("N109").Value = -11979.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7715707
$ws.Range("I11").Value = 10000490
$ws.Range("J11").Value = 2003750
$ws.Range("K11").Value = 10000490
$ws.Range("L11").Value = 2003750
$ws.Range("M11").Value = -10000351

$ws.Range("H70").Value = 5004.9243
$ws.Range("I70").Value = 4411.4546
$ws.Range("J70").Value = 5984.15
$ws.Range("K70").Value = 4411.4546
$ws.Range("L70").Value = 5984.15
$ws.Range("M70").Value = -4141.4546
$ws.Range("N70").Value = -6524.15

$ws.Range("H73").Value = 5004.9243
$ws.Range("I73").Value = 4411.4546
$ws.Range("J73").Value = 5984.15
$ws.Range("K73").Value = 4411.4546
$ws.Range("L73").Value = 5984.15
$ws.Range("M73").Value = -3475.4546
$ws.Range("N73").Value = -7856.15

$ws.Range("H80").Value = 4278.2144
$ws.Range("I80").Value = 5165.8335
$ws.Range("J80").Value = 3612.5
$ws.Range("K80").Value = 5165.8335
$ws.Range("L80").Value = 3612.5
$ws.Range("M80").Value = -4167.8335
$ws.Range("N80").Value = -5608.5

$ws.Range("H83").Value = 4278.2144
$ws.Range("I83").Value = 5165.8335
$ws.Range("J83").Value = 3612.5
$ws.Range("K83").Value = 25829.1675
$ws.Range("L83").Value = 18062.5
$ws.Range("M83").Value = -20837.1675
$ws.Range("N83").Value = -28046.5

$ws.Range("H113").Value = 1801.6666
$ws.Range("I113").Value = 1741.5834
$ws.Range("J113").Value = 2042
$ws.Range("K113").Value = 1741.5834
$ws.Range("L113").Value = 2042
$ws.Range("M113").Value = 428.4166
$ws.Range("N113").Value = -6382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 77745.69500000001
$ws.Range("I46").Value = 143513.14
$ws.Range("J46").Value = 1017
$ws.Range("K46").Value = 143513.14
$ws.Range("L46").Value = 1017
$ws.Range("M46").Value = -143325.14
$ws.Range("N46").Value = -1393

$ws.Range("H61").Value = 2325
$ws.Range("I61").Value = 2100
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2100
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1898
$ws.Range("N61").Value = -3404

$ws.Range("H82").Value = 2554.923
$ws.Range("I82").Value = 958.7143
$ws.Range("J82").Value = 4417.1665
$ws.Range("K82").Value = 958.7143
$ws.Range("L82").Value = 4417.1665
$ws.Range("M82").Value = -597.7143
$ws.Range("N82").Value = -5139.1665

$ws.Range("H85").Value = 2554.923
$ws.Range("I85").Value = 958.7143
$ws.Range("J85").Value = 4417.1665
$ws.Range("K85").Value = 958.7143
$ws.Range("L85").Value = 4417.1665
$ws.Range("M85").Value = 289.2857
$ws.Range("N85").Value = -6913.1665

$ws.Range("H113").Value = 2325
$ws.Range("I113").Value = 2100
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
